$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated data rows (R007/R008 placeholder rows replaced with real report
# rows, plus two additional rows that were missing from the export).
$data = @(
    @("R001", "Jane Doe", "Graduate School", "Master", "Complete", "2023-05-04 23:21:13", 1),
    @("R002", "John Doe", "Senior High School", "Accountancy, Business, and Management Strand", "Complete", "2023-05-04 23:21:34", 1),
    @("R003", "Mark Doe", "College", "Bachelor of Science in Information Technology", "Complete", "2023-05-04 23:21:53", 1),
    @("R004", "Son Doe", "Junior High School", "Junior High School", "Complete", "2023-05-04 23:22:11", 1)
)

$row = 2
foreach ($record in $data) {
    $ws.Cells.Item($row, 1).Value = $record[0]
    $ws.Cells.Item($row, 2).Value = $record[1]
    $ws.Cells.Item($row, 3).Value = $record[2]
    $ws.Cells.Item($row, 4).Value = $record[3]
    $ws.Cells.Item($row, 5).Value = $record[4]
    $ws.Cells.Item($row, 6).Value = $record[5]
    $ws.Cells.Item($row, 7).Value = $record[6]
    $row++
}
